$wb = $excel.ActiveWorkbook

# --- Sheet "Country": update row 2 values (B2:E2) ---
$wsCountry = $wb.Worksheets.Item("Country")
$wsCountry.Cells.Item(2, 2).Value = 21.09
$wsCountry.Cells.Item(2, 3).Value = 19.73
$wsCountry.Cells.Item(2, 4).Value = 13.2
$wsCountry.Cells.Item(2, 5).Value = 9.66

# --- Sheet "States": rewrite state rows in new order/values, one fewer row (Ladakh removed) ---
$wsStates = $wb.Worksheets.Item("States")

$statesData = @(
    @("IN-BR", "Bihar", 0, 7.89, 7.89, 39.47, 44.74),
    @("IN-JH", "Jharkhand", 0, 4.17, 29.17, 33.33, 29.17),
    @("IN-AR", "Arunachal Pradesh", 0, 24, 16, 4, 28),
    @("IN-MZ", "Mizoram", 0, 27.27, 27.27, 18.18, 27.27),
    @("IN-UP", "Uttar Pradesh", 0, 10.67, 22.67, 33.33, 25.33),
    @("IN-PB", "Punjab", 0, 18.18, 36.36, 13.64, 22.73),
    @("IN-MN", "Manipur", 0, 18.75, 31.25, 31.25, 18.75),
    @("IN-NL", "Nagaland", 0, 27.27, 18.18, 27.27, 9.09),
    @("IN-ML", "Meghalaya", 0, 27.27, 36.36, 18.18, 9.09),
    @("IN-TS", "Telangana", 0, 30.3, 15.15, 18.18, 9.09),
    @("IN-JK", "Jammu and Kashmir", 0, 31.82, 36.36, 13.64, 9.09),
    @("IN-OR", "Odisha", 0, 36.67, 30, 10, 6.67),
    @("IN-MP", "Madhya Pradesh", 0, 34.62, 28.85, 11.54, 1.92),
    @("IN-DD", "Daman and Diu", 0, 0, 50, 50, 0),
    @("IN-PY", "Puducherry", 0, 0, 25, 50, 0),
    @("IN-AS", "Assam", 0, 30.3, 18.18, 15.15, 0),
    @("IN-TR", "Tripura", 0, 25, 62.5, 12.5, 0),
    @("IN-HR", "Haryana", 0, 18.18, 54.55, 9.09, 0),
    @("IN-DL", "Delhi", 0, 36.36, 45.45, 9.09, 0),
    @("IN-HP", "Himachal Pradesh", 0, 25, 16.67, 8.33, 0),
    @("IN-WB", "West Bengal", 0, 30.43, 8.699999999999999, 4.35, 0),
    @("IN-CT", "Chhattisgarh", 0, 44.44, 22.22, 3.7, 0),
    @("IN-RJ", "Rajasthan", 0, 33.33, 21.21, 0, 0),
    @("IN-MH", "Maharashtra", 0, 19.44, 11.11, 0, 0),
    @("IN-GJ", "Gujarat", 0, 12.12, 6.06, 0, 0),
    @("IN-KA", "Karnataka", 0, 6.67, 3.33, 0, 0),
    @("IN-TN", "Tamil Nadu", 0, 10.81, 2.7, 0, 0),
    @("IN-CH", "Chandigarh", 0, 100, 0, 0, 0),
    @("IN-DN", "Dadra and Nagar Haveli", 0, 100, 0, 0, 0),
    @("IN-AP", "Andhra Pradesh", 0, 15.38, 0, 0, 0),
    @("IN-UL", "Uttarakhand", 0, 7.69, 0, 0, 0)
)

for ($i = 0; $i -lt $statesData.Count; $i++) {
    $rowNum = $i + 2
    $rowVals = $statesData[$i]
    $wsStates.Cells.Item($rowNum, 1).Value = $rowVals[0]
    $wsStates.Cells.Item($rowNum, 2).Value = $rowVals[1]
    $wsStates.Cells.Item($rowNum, 3).Value = $rowVals[2]
    $wsStates.Cells.Item($rowNum, 4).Value = $rowVals[3]
    $wsStates.Cells.Item($rowNum, 5).Value = $rowVals[4]
    $wsStates.Cells.Item($rowNum, 6).Value = $rowVals[5]
    $wsStates.Cells.Item($rowNum, 7).Value = $rowVals[6]
}

# Remove the now-extra last row (previously row 33, Uttarakhand duplicate/Ladakh dropped)
$wsStates.Rows.Item(33).Delete()

# --- Sheet "Dark clusters": update rows 2-4 ---
$wsDark = $wb.Worksheets.Item("Dark clusters")
$wsDark.Cells.Item(2, 3).Value = 4.44
$wsDark.Cells.Item(2, 4).Value = 8.89
$wsDark.Cells.Item(2, 5).Value = 33.33
$wsDark.Cells.Item(2, 6).Value = 53.33

$wsDark.Cells.Item(3, 4).Value = 60
$wsDark.Cells.Item(3, 5).Value = 15
$wsDark.Cells.Item(3, 6).Value = 25

$wsDark.Cells.Item(4, 3).Value = 26.09
$wsDark.Cells.Item(4, 4).Value = 25.22
$wsDark.Cells.Item(4, 5).Value = 16.52
$wsDark.Cells.Item(4, 6).Value = 13.04

